# refs #882 Video Wall -> HSR Videowall
# Rename the "Raw data" table's header columns (typing directly into the
# header cells, the way a user would in Excel - this keeps the backing
# ListObject / table definition in sync automatically).

$wb = $excel.ActiveWorkbook
$wsRaw = $wb.Worksheets.Item("Raw data")

$wsRaw.Range("A3").Value = "Anzahl Monitore"
$wsRaw.Range("B3").Value = "Treibermodell"
$wsRaw.Range("C3").Value = "Monitormodus"
$wsRaw.Range("D3").Value = "Videogrösse" + [char]10 + "(x*FullHD)"

# Leave the final selection where the user ended up after editing the
# headers on the "Raw data" sheet.
$wsRaw.Range("E13").Select()
